# Cookie Checker Upload Compare and Send Email
# Update the cookie-check results for three domains:
#  - https://www.nbg.gr/go4more/                              (row 3)  -> cookie count 1 -> 2, add "WSS_FullScreenMode"
#  - https://www.nbg.gr/el/retail/housing-loans/Calculator/    (row 7)  -> cookie count 1 -> 2, add "NBGPublicSite"
#  - https://microsites.nbg.gr/DonationApply                   (row 30) -> cookie count 0 -> 1, add "NBGPUBLICConsent"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: https://www.nbg.gr/go4more/
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = "WSS_FullScreenMode"

# Row 7: https://www.nbg.gr/el/retail/housing-loans/Calculator/
$ws.Range("B7").Value = 2
$ws.Range("D7").Value = "NBGPublicSite"

# Row 30: https://microsites.nbg.gr/DonationApply
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "NBGPUBLICConsent"
